$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" sheet, which both list the same events.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 288
    $ws.Range("F5").Value = 4236
}
